$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1. Split/rewrite the "Model" section paragraph into three paragraphs
# ---------------------------------------------------------------
$rng = $d.Content
$old = ", Manager and the Clients. We will model them as actors "
$new = " and the Clients. Each of them will have a state associated to them that describes what they are doing in each moment which will influence how they interact with each other. For example: a waiter that is busy cleaning a table may not serve a drink prepared by the barman. If the waiter was idle in the home he would immediately take the drink from the barman and bring it to the client when he receives the request.`rThe various entities of the system will be coordinated by exchanging messages of various kinds (dispatches, request/reply and events). These messages will let us model various activities such as informing the client of the maximum waiting time, taking the order of the client and transmitting it to the barman.`rWe will model the entities as Actors "
$found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

# Remove the now-duplicated old paragraph (originally held "The various entities..."
# text) plus the two now-empty paragraphs that used to follow it.
$start = $d.Paragraphs.Item(28).Range.Start
$end = $d.Paragraphs.Item(30).Range.End
$d.Range($start, $end).Delete()

for ($i = 22; $i -le 28; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Host $i ": [" $p.Range.Text "]"
}

# ---------------------------------------------------------------
# 2. Underline "and" inside paragraph 26 (the "various entities..." one)
# ---------------------------------------------------------------
$andRng = $d.Content
$found = $andRng.Find.Execute("reply and events", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Found reply/events:" $found "[" $andRng.Text "]" $andRng.Start $andRng.End
# "reply and events" -> "and" starts right after "reply " (6 chars)
$andStart = $andRng.Start + 6
$andEnd = $andStart + 3
$andOnly = $d.Range($andStart, $andEnd)
Write-Host "and-only range: [" $andOnly.Text "]"
$andOnly.Font.Underline = 1

# ---------------------------------------------------------------
# 3. Underline "entities" and "Actors" inside paragraph 27 ("We will model...")
#    and give the paragraph mark itself the underline property too.
# ---------------------------------------------------------------
$p27 = $d.Paragraphs.Item(27)
Write-Host "p27 text: [" $p27.Range.Text "]"

$entitiesRng = $d.Content
$f3 = $entitiesRng.Find.Execute("entities as Actors", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Found entities-as-Actors:" $f3 "[" $entitiesRng.Text "]" $entitiesRng.Start $entitiesRng.End

$entitiesStart = $entitiesRng.Start
$entitiesEnd = $entitiesStart + 8   # "entities" is 8 chars
$entitiesOnly = $d.Range($entitiesStart, $entitiesEnd)
Write-Host "entities-only: [" $entitiesOnly.Text "]"
$entitiesOnly.Font.Underline = 1

$actorsStart = $entitiesRng.End - 6  # "Actors" is 6 chars
$actorsEnd = $entitiesRng.End
$actorsOnly = $d.Range($actorsStart, $actorsEnd)
Write-Host "actors-only: [" $actorsOnly.Text "]"
$actorsOnly.Font.Underline = 1

# Give the paragraph MARK itself the underline property (matches <w:pPr><w:rPr><w:u .../>)
# without underlining the rest of the paragraph's visible text: insert a
# temporary character right at the end of the paragraph, underline it
# together with the mark, then delete just that temporary character again
# (Word keeps the formatting on the orphaned paragraph mark).
$p27 = $d.Paragraphs.Item(27)
$endRng = $p27.Range
$endRng.InsertBefore("Z")
$markRng = $d.Range($endRng.Start, $endRng.End)
Write-Host "markRng (with Z) text: [" $markRng.Text "]"
$charOnly = $d.Range($markRng.Start, $markRng.Start + 1)
Write-Host "char only: [" $charOnly.Text "]"
$charOnly.Delete()

